$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17; existing rows 17-61 shift down to 18-62.
$ws.Rows(17).Insert()

# Populate the newly inserted row 17 with the new record.
$ws.Range("A17").Value = 10
$ws.Range("B17").Value = "Vega Modelo de Temuco"
$ws.Range("C17").Value = "La Araucanía"
$ws.Range("D17").Value = 44792
$ws.Range("E17").Value = 9
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100108
$ws.Range("H17").Value = "Tropicales y subtropicales"
$ws.Range("I17").Value = 100108003
$ws.Range("J17").Value = "Maracuyá"
$ws.Range("K17").Value = "Sin especificar"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 12
$ws.Range("N17").Value = 36000
$ws.Range("O17").Value = 36000
$ws.Range("P17").Value = 36000
$ws.Range("Q17").Value = "$/caja 18 kilos"
$ws.Range("R17").Value = "Región de Arica y Parinacota"
$ws.Range("S17").Value = 2000
$ws.Range("T17").Value = 18
